# ---------------------------------------------------------------------------
# NC92Soil 0.8 (29/09/2020)
#  - stochastic batch-input generation sheet ("Stochastic")
#  - list-of-valid-values helper sheet ("Voices")
#  - "Profiles" sheet populated with multi-profile definitions
#  - "Clusters" sheet: obsolete IW1/SUB1 sample row removed
#  - "Soils" sheet: new "Curve Std" column
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheets: add "Stochastic" in front, add "Voices" at the end.
# ---------------------------------------------------------------------------
$soils = $wb.Worksheets.Item("Soils")
$stoch = $wb.Worksheets.Add($soils)
$stoch.Name = "Stochastic"

$voices = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$voices.Name = "Voices"

# ---------------------------------------------------------------------------
# 2. "Voices" sheet - simple list used by the data validation below.
# ---------------------------------------------------------------------------
$voices.Range("A1").Value = "All profile"
$voices.Range("A2").Value = "Single groups"
$voices.Columns.Item(1).ColumnWidth = 11.55

# ---------------------------------------------------------------------------
# 3. "Stochastic" sheet - headers, sample data, styling, comments and the
#    data-validation list pointing at the "Voices" sheet.
# ---------------------------------------------------------------------------

# -- headers (row 1) --------------------------------------------------------
$stoch.Range("A1").Value = "Group name"
$stoch.Range("B1").Value = "Unit weight`n[KN/m3]"
$stoch.Range("C1").Value = "Min thickness`n[m]"
$stoch.Range("D1").Value = "Max thickness`n[m]"
$stoch.Range("E1").Value = "Vs Law"
$stoch.Range("F1").Value = "Vs Std"
$stoch.Range("G1").Value = "Degradation curve`nMean"
$stoch.Range("H1").Value = "Degradation curve`nStd"
$stoch.Range("I1").Value = "Inter-layer correlation"
$stoch.Range("M1").Value = "Number of iterations"
$stoch.Range("N1").Value = "Random seed"
$stoch.Range("O1").Value = "Correlation mode"

# -- header styling: centred, vertically centred; most of them wrap --------
foreach ($col in @("A", "F")) {
    $r = $stoch.Range($col + "1")
    $r.VerticalAlignment = -4108
    $r.HorizontalAlignment = -4108
}
foreach ($col in @("B", "C", "D", "E", "G", "H", "I")) {
    $r = $stoch.Range($col + "1")
    $r.VerticalAlignment = -4108
    $r.HorizontalAlignment = -4108
    $r.WrapText = $true
}
$stoch.Rows.Item(1).RowHeight = 43.2

# -- data rows ---------------------------------------------------------------
$stoch.Range("A2").Value = "A"
$stoch.Range("B2").Value = 18
$stoch.Range("C2").Value = 5
$stoch.Range("D2").Value = 10
$stoch.Range("E2").Value = "10*x + 100"
$stoch.Range("F2").Value = 50
$stoch.Range("G2").Value = "Vucetic & Dobry, PI=0"
$stoch.Range("H2").Value = 1
$stoch.Range("I2").Value = "Toro: USGS AB"
$stoch.Range("M2").Value = 100
$stoch.Range("O2").Value = "Single groups"

$stoch.Range("A3").Value = "G"
$stoch.Range("B3").Value = 19
$stoch.Range("C3").Value = 5
$stoch.Range("D3").Value = 10
$stoch.Range("E3").Value = "0.45*x^2 + 200"
$stoch.Range("F3").Value = 70
$stoch.Range("G3").Value = "Vucetic & Dobry, PI=0"
$stoch.Range("H3").Value = 1
$stoch.Range("I3").Value = "Toro: USGS AB"

$stoch.Range("A4").Value = "S"
$stoch.Range("B4").Value = 18
$stoch.Range("C4").Value = 5
$stoch.Range("D4").Value = 10
$stoch.Range("E4").Value = "1.35^x + 250"
$stoch.Range("F4").Value = 90
$stoch.Range("G4").Value = "Vucetic & Dobry, PI=0"
$stoch.Range("H4").Value = 1
$stoch.Range("I4").Value = "Toro: USGS AB"

# -- "Vs Law" / "Inter-layer correlation" values are right-aligned ----------
$stoch.Range("E2:E4").VerticalAlignment = -4108
$stoch.Range("E2:E4").HorizontalAlignment = -4152
$stoch.Range("I2:I4").VerticalAlignment = -4108
$stoch.Range("I2:I4").HorizontalAlignment = -4152

# -- column widths ------------------------------------------------------------
$stoch.Columns.Item("A").ColumnWidth = 16.66
$stoch.Columns.Item("C").ColumnWidth = 15.11
$stoch.Columns.Item("D").ColumnWidth = 14.11
$stoch.Columns.Item("E").ColumnWidth = 16.11
$stoch.Columns.Item("F").ColumnWidth = 11.22
$stoch.Columns.Item("G").ColumnWidth = 18.78
$stoch.Columns.Item("H").ColumnWidth = 16.55
$stoch.Columns.Item("I").ColumnWidth = 16.55
$stoch.Columns.Item("K").ColumnWidth = 5.56
$stoch.Columns.Item("L").ColumnWidth = 5.11
$stoch.Columns.Item("M").ColumnWidth = 18.22
$stoch.Columns.Item("N").ColumnWidth = 12
$stoch.Columns.Item("O").ColumnWidth = 15.44

# -- cell comments (instructions) -------------------------------------------
$wb.Author = "Gianluca Acunzo"

$stoch.Range("E1").AddComment("Gianluca Acunzo:`n`nSet the law which relates the layer depth with its mean Vs. Use ""x"" as the variable for depth`nEx.`n3*x^2 + 4*x + 9") | Out-Null

$stoch.Range("I1").AddComment("Gianluca Acunzo:`n`nThe inter-layer correlation law. Write an expression using ""x"" for the depth and ""y"" for the thickness of the layer. `nIf ""Toro: [model]"" is specified, the Toro velocity variation model is used according to the parameters of [model]. The mean value of Vs and its standard deviation are obtained from given data and not from default model.") | Out-Null

$stoch.Range("O1").AddComment("Gianluca Acunzo:`n`nIf ""All profile"" is selected, the inter-layer correlation specified in first row will be used for the whole profile.`nIf ""Single groups"" is selected, the correlations will be considered independently for each group.") | Out-Null

# -- data validation: O2 restricted to the "Voices" sheet list --------------
$stoch.Range("O2").Validation.Add(3, 1, 1, "=Voices!`$A`$1:`$A`$2")

$stoch.Range("I10").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. "Soils" sheet - add the "Curve Std" column header.
# ---------------------------------------------------------------------------
$soils.Range("G1").Value = "Curve Std"
$soils.Range("G1").VerticalAlignment = -4108
$soils.Range("G1").HorizontalAlignment = -4108
$soils.Range("G1").WrapText = $true

# ---------------------------------------------------------------------------
# 5. "Clusters" sheet - drop the obsolete IW1/SUB1 sample row.
# ---------------------------------------------------------------------------
$clusters = $wb.Worksheets.Item("Clusters")
$clusters.Rows.Item(2).Delete()

# ---------------------------------------------------------------------------
# 6. "Profiles" sheet - populate with the P1/P2/P3 profile definitions.
# ---------------------------------------------------------------------------
$profiles = $wb.Worksheets.Item("Profiles")

$profiles.Range("A1").Value = "P1"
$profiles.Range("B1").Value = "P2"
$profiles.Range("C1").Value = "P3"
$profiles.Range("A1:C1").Font.Bold = $true
$profiles.Range("A1:C1").VerticalAlignment = -4108
$profiles.Range("A1:C1").HorizontalAlignment = -4108

$profiles.Range("A2").Value = "Spettro UHS 2.txt"
$profiles.Range("B2").Value = "Spettro UHS 2.txt"
$profiles.Range("C2").Value = "Spettro UHS 2.txt; Spettro UHS 3.txt"

$profiles.Range("A3").Value = "A;9;250"
$profiles.Range("B3").Value = "A;5"
$profiles.Range("C3").Value = "S;3"

$profiles.Range("A4").Value = "G;8;300"
$profiles.Range("B4").Value = "S;9"
$profiles.Range("C4").Value = "A;5"

$profiles.Range("A5").Value = "A;3;350"
$profiles.Range("B5").Value = "A;8"
$profiles.Range("C5").Value = "G;9"

# ---------------------------------------------------------------------------
# 7. Final selections / active sheet, matching the authored workbook.
# ---------------------------------------------------------------------------
$soils.Activate()
$wb.Worksheets.Item("Clusters").Range("G4").Select() | Out-Null
$soils.Select() | Out-Null
